# ---------------------------------------------------------------------------
# "merged with spark jobs for embeddings"
#
# Inserts a new "Word Weighting" column (D) ahead of the metric columns,
# re-creates the header hyperlinks that shifted one column to the right,
# refreshes the HR@10 / nDCG@10 / CSHR@10 / CSnDCG@10 numbers for the
# existing LSI / Prod2Vec rows and appends the new ProdWord2Vec results
# (uniform + tfidf word-weighting, rows 11-16).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Make room for the new "Word Weighting" column. This shifts the old
#    D,E,F,G (HR@10, nDCG@10, CSHR@10, CSnDCG@10) columns -> E,F,G,H,
#    including their values/styles, but NOT the hyperlinks (Excel's COM
#    model does not re-target hyperlinks on an Insert, so we rebuild them
#    below).
$ws.Columns("D").Insert() | Out-Null

# 2. New column header + the four metric headers that now live one column
#    further right.
$ws.Range("D1").Value = "Word Weighting"

# 3. Re-create the header hyperlinks against their new home cells, reusing
#    the exact same mailto targets as before.
$ws.Hyperlinks.Delete() | Out-Null
$ws.Hyperlinks.Add($ws.Range("E1"), "mailto:HR@10") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F1"), "mailto:nDCG@10") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G1"), "mailto:CSHR@10") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H1"), "mailto:CSnDCG@10") | Out-Null

# 4. Refreshed metric values for the existing LSI (rows 2-4) and Prod2Vec
#    (rows 5-10) results.
$ws.Range("E2").Value = 0.00100038476337052
$ws.Range("F2").Value = 0.0056671033625361397
$ws.Range("G2").Value = 0.00100038476337052
$ws.Range("H2").Value = 0.0056671033625361397

$ws.Range("E3").Value = 0.00159676798768756
$ws.Range("F3").Value = 0.0088026570135895394
$ws.Range("G3").Value = 0.0015582916506348501
$ws.Range("H3").Value = 0.0086785184961932493

$ws.Range("E4").Value = 0.0013466717968449299
$ws.Range("F4").Value = 0.0082923226745443002
$ws.Range("G4").Value = 0.0013466717968449299
$ws.Range("H4").Value = 0.0082923226745443002

$ws.Range("E5").Value = 0.0107541362062332
$ws.Range("F5").Value = 0.049252926981654101
$ws.Range("G5").Value = 0.00107733743747595
$ws.Range("H5").Value = 0.0058850580428803896

$ws.Range("E6").Value = 0.0126394767218161
$ws.Range("F6").Value = 0.060849386077138197
$ws.Range("G6").Value = 0.00080800307810696399
$ws.Range("H6").Value = 0.0043796206328408803

$ws.Range("E7").Value = 0.012524047710657999
$ws.Range("F7").Value = 0.065447976587353504
$ws.Range("G7").Value = 0.000865717583686033
$ws.Range("H7").Value = 0.0047150872562152996

$ws.Range("E8").Value = 0.0098114659484417905
$ws.Range("F8").Value = 0.046397192349962403
$ws.Range("G8").Value = 0.00069257406694882599
$ws.Range("H8").Value = 0.0046865417366717304

$ws.Range("E9").Value = 0.014024624855713701
$ws.Range("F9").Value = 0.067112001487567199
$ws.Range("G9").Value = 0.00111581377452866
$ws.Range("H9").Value = 0.0056959708234617301

$ws.Range("E10").Value = 0.013639861485186601
$ws.Range("F10").Value = 0.063758922736521395
$ws.Range("G10").Value = 0.000848133863686033
$ws.Range("H10").Value = 0.0047079185862153004

# 5. New ProdWord2Vec rows (Spark-produced embeddings), word weighting is
#    recorded per-row: uniform weighting first (vector sizes 50/100/200),
#    then tf-idf weighting (vector sizes 50/100/200) - this mirrors the
#    order the source data was appended in, so shared-string indices line
#    up with the author's original save.
$ws.Range("A11").Value = "ProdWord2Vec"
$ws.Range("B11").Value = 50
$ws.Range("D11").Value = "uniform"
$ws.Range("E11").Value = 0.00846479415159683
$ws.Range("F11").Value = 0.048245424301151703
$ws.Range("G11").Value = 0.0023470565602154601
$ws.Range("H11").Value = 0.011803630366097001

$ws.Range("A12").Value = "ProdWord2Vec"
$ws.Range("B12").Value = 100
$ws.Range("D12").Value = "uniform"
$ws.Range("E12").Value = 0.0083493651404386894
$ws.Range("F12").Value = 0.050031575072219003
$ws.Range("G12").Value = 0.0027125817622162302
$ws.Range("H12").Value = 0.0136490289779198

$ws.Range("A13").Value = "ProdWord2Vec"
$ws.Range("B13").Value = 200
$ws.Range("D13").Value = "uniform"
$ws.Range("E13").Value = 0.0077145055790689299
$ws.Range("F13").Value = 0.045540062899378703
$ws.Range("G13").Value = 0.0028664871104270798
$ws.Range("H13").Value = 0.013932971840025499

$ws.Range("A14").Value = "ProdWord2Vec"
$ws.Range("B14").Value = 50
$ws.Range("D14").Value = "tfidf"
$ws.Range("E14").Value = 0.0096383224317045796
$ws.Range("F14").Value = 0.054303058820772
$ws.Range("G14").Value = 0.0026741054251635098
$ws.Range("H14").Value = 0.0123542901107217

$ws.Range("A15").Value = "ProdWord2Vec"
$ws.Range("B15").Value = 100
$ws.Range("D15").Value = "tfidf"
$ws.Range("E15").Value = 0.0095806079261255093
$ws.Range("F15").Value = 0.054419743212762203
$ws.Range("G15").Value = 0.0033089649865332702
$ws.Range("H15").Value = 0.016322618366025901

$ws.Range("A16").Value = "ProdWord2Vec"
$ws.Range("B16").Value = 200
$ws.Range("D16").Value = "tfidf"
$ws.Range("E16").Value = 0.0097537514428627202
$ws.Range("F16").Value = 0.054268502405740902
$ws.Range("G16").Value = 0.0031550596383224201
$ws.Range("H16").Value = 0.016580330222818799

# 6. Bold the best-performing cell per metric, per the source notebook's
#    manual "highlight the winner" pass.
$ws.Range("G3,H3,E9,F9").Font.Bold = $true
$ws.Range("E15,F15,G15,H15").Font.Bold = $true

# 7. Match the print orientation that was set alongside this data refresh.
$ws.PageSetup.Orientation = 1

# 8. Leave the same cell selected as in the saved workbook.
$ws.Range("H15").Select() | Out-Null
